$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("PE, SEPTEMBER")

# Update row 19
$ws.Range("C19").Value = 45927
$ws.Range("G19").Value = 517948038
$ws.Range("I19").Formula = "=1353132-56255.04"

# Update row 20
$ws.Range("C20").Value = 45927
$ws.Range("G20").Value = 517947574
$ws.Range("I20").Formula = "=1353132-56255.04"

# Update row 21
$ws.Range("C21").Value = 45927
$ws.Range("G21").Value = 517947435
$ws.Range("I21").Formula = "=1353132-56255.04"

# Update selection to H21
$ws.Range("H21").Select()
